# Call_CMT.xlsx update: "add functionality for rsm"
#
# The CMT breakdown table (Sheet1) is refreshed with new data:
#   - two additional CMT categories (CMT14 and CMT36) are introduced, which
#     pushes the table from 32 data rows (rows 2-33) to 34 data rows
#     (rows 2-35), and
#   - every numeric cell in the B:J columns is refreshed with new figures.
#
# This script rewrites the whole data block (header stays as-is) with the
# final values, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is: row number, CMT label (column A), then the 9 numeric
# values for columns B..J in order.
$data = @(
    @(2, "CMT11", 1, 0, 0, 0, 0, 1, 1, 0, 0),
    @(3, "CMT12", 3, 10, 6, 3, 8, 1, 10, 5, 6),
    @(4, "CMT13", 15, 14, 15, 2, 14, 14, 15, 8, 12),
    @(5, "CMT14", 0, 0, 0, 1, 0, 0, 0, 2, 0),
    @(6, "CMT15", 8, 4, 0, 3, 5, 7, 10, 0, 2),
    @(7, "CMT16", 8, 2, 0, 7, 7, 6, 9, 1, 1),
    @(8, "CMT21", 3, 6, 0, 0, 8, 2, 4, 3, 3),
    @(9, "CMT22", 10, 3, 7, 4, 3, 16, 8, 5, 7),
    @(10, "CMT23", 11, 6, 9, 2, 0, 5, 0, 5, 5),
    @(11, "CMT24", 8, 0, 5, 1, 5, 6, 3, 5, 3),
    @(12, "CMT25", 6, 3, 4, 1, 4, 4, 11, 1, 6),
    @(13, "CMT26", 6, 5, 7, 6, 4, 7, 3, 17, 0),
    @(14, "CMT31", 1, 7, 3, 7, 11, 8, 4, 11, 0),
    @(15, "CMT32", 0, 4, 0, 0, 10, 5, 5, 0, 0),
    @(16, "CMT33", 2, 0, 7, 6, 4, 5, 3, 0, 5),
    @(17, "CMT34", 8, 7, 10, 0, 3, 12, 3, 0, 14),
    @(18, "CMT35", 0, 0, 0, 0, 0, 0, 0, 0, 1),
    @(19, "CMT36", 3, 6, 4, 1, 5, 8, 6, 0, 5),
    @(20, "CMT41", 10, 10, 8, 3, 5, 3, 8, 1, 5),
    @(21, "CMT42", 12, 14, 8, 17, 12, 20, 14, 6, 4),
    @(22, "CMT43", 4, 5, 7, 7, 13, 5, 0, 3, 6),
    @(23, "CMT44", 6, 5, 7, 3, 11, 2, 9, 5, 4),
    @(24, "CMT45", 8, 7, 7, 4, 4, 6, 4, 3, 10),
    @(25, "CMT46", 6, 7, 3, 5, 13, 16, 5, 11, 12),
    @(26, "CMT51", 3, 4, 2, 9, 11, 12, 11, 5, 0),
    @(27, "CMT52", 5, 7, 7, 1, 6, 8, 10, 2, 5),
    @(28, "CMT53", 8, 5, 4, 1, 8, 3, 6, 0, 3),
    @(29, "CMT54", 6, 11, 1, 5, 0, 13, 12, 0, 2),
    @(30, "CMT55", 4, 7, 4, 2, 5, 4, 10, 3, 6),
    @(31, "CMT61", 3, 5, 2, 0, 7, 6, 3, 4, 4),
    @(32, "CMT62", 9, 4, 10, 6, 12, 7, 9, 13, 8),
    @(33, "CMT63", 3, 0, 1, 3, 2, 4, 2, 0, 0),
    @(34, "CMT64", 10, 10, 5, 10, 10, 7, 5, 10, 2),
    @(35, "CMT65", 9, 6, 6, 4, 4, 11, 12, 6, 7)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}
